$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cyclic rotation of weekly data between rows 2, 3 and 4:
# new row2 = old row3, new row3 = old row4, new row4 = old row2
# (columns D, J, K, L, M, P)

$ws.Range("D2").Value = 44827
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 30000
$ws.Range("L2").Value = 31000
$ws.Range("M2").Value = 30500
$ws.Range("P2").Value = 1220

$ws.Range("D3").Value = 44414
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 31000
$ws.Range("L3").Value = 32000
$ws.Range("M3").Value = 31500
$ws.Range("P3").Value = 1260

$ws.Range("D4").Value = 44379
$ws.Range("J4").Value = 240
$ws.Range("K4").Value = 31000
$ws.Range("L4").Value = 32000
$ws.Range("M4").Value = 31500
$ws.Range("P4").Value = 1260
